$d = $word.ActiveDocument

# The "Layer building" table is the second table in the document. It has
# a header row ("NO" / "Name of column" / "Description") followed by six
# data rows numbered 1-6. Rows 5 ("BCR" / Maximum construction density)
# and 6 ("FAR" / Coefficient of land use) are being removed, leaving the
# table ending right after row 4 ("Height" / Maximum height in the
# planning (meter)).

$t = $d.Tables.Item(2)

# Row 6 (index 7 overall, counting the header row) is "FAR" / Coefficient
# of land use in the planning - delete it first so row indices for the
# earlier "BCR" row remain stable.
$t.Rows.Item(7).Delete()

# Row 5 (index 6 overall) is "BCR" / Maximum construction density in the
# planning (%).
$t.Rows.Item(6).Delete()
